# Apply Chocobo_Profits market-data refresh (scheduled runner update)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 32475
$ws.Range("J81").Value = 32475
$ws.Range("L81").Value = 32475
$ws.Range("N81").Value = -34471

$ws.Range("H84").Value = 32475
$ws.Range("J84").Value = 32475
$ws.Range("L84").Value = 97425
$ws.Range("N84").Value = -107409

$ws.Range("H93").Value = 23741.234
$ws.Range("J93").Value = 23741.234
$ws.Range("L93").Value = 23741.234
$ws.Range("N93").Value = -28733.234

$ws.Range("H107").Value = 3603.8333
$ws.Range("I107").Value = 4822.5
$ws.Range("J107").Value = 1166.5
$ws.Range("K107").Value = 4822.5
$ws.Range("L107").Value = 1166.5
$ws.Range("M107").Value = -2902.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4213.9443
$ws.Range("I32").Value = 4154.7085
$ws.Range("J32").Value = 4687.8335
$ws.Range("K32").Value = 4154.7085
$ws.Range("L32").Value = 4687.8335
$ws.Range("M32").Value = -3867.7085
$ws.Range("N32").Value = -5261.8335

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H103").Value = 33726.77
$ws.Range("J103").Value = 33726.77
$ws.Range("L103").Value = 33726.77
$ws.Range("N103").Value = -36070.77

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 2554.2273
$ws.Range("J132").Value = 3435.3572
$ws.Range("L132").Value = 10306.0716
$ws.Range("N132").Value = -15366.0716

$ws.Range("H137").Value = 40001.668
$ws.Range("J137").Value = 40001.668
$ws.Range("L137").Value = 40001.668
$ws.Range("N137").Value = -50201.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 19000
$ws.Range("J76").Value = 19000
$ws.Range("L76").Value = 19000
$ws.Range("N76").Value = -19630

$ws.Range("H79").Value = 19000
$ws.Range("J79").Value = 19000
$ws.Range("L79").Value = 19000
$ws.Range("N79").Value = -21184

$ws.Range("H92").Value = 71249.75
$ws.Range("J92").Value = 71249.75
$ws.Range("L92").Value = 71249.75
$ws.Range("N92").Value = -76241.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 823.8461
$ws.Range("I22").Value = 202.2
$ws.Range("J22").Value = 1212.375
$ws.Range("K22").Value = 202.2
$ws.Range("L22").Value = 1212.375
$ws.Range("M22").Value = 147.8
$ws.Range("N22").Value = -1912.375

$ws.Range("H68").Value = 85073
$ws.Range("J68").Value = 85073
$ws.Range("L68").Value = 85073
$ws.Range("N68").Value = -86571

$ws.Range("H71").Value = 85073
$ws.Range("J71").Value = 85073
$ws.Range("L71").Value = 255219
$ws.Range("N71").Value = -262707

$ws.Range("H87").Value = 23214.285
$ws.Range("J87").Value = 23214.285
$ws.Range("L87").Value = 23214.285
$ws.Range("N87").Value = -25586.285

$ws.Range("H90").Value = 23214.285
$ws.Range("J90").Value = 23214.285
$ws.Range("L90").Value = 69642.855
$ws.Range("N90").Value = -81498.855

$ws.Range("H99").Value = 15389792
$ws.Range("I99").Value = 40002380
$ws.Range("J99").Value = 6925
$ws.Range("K99").Value = 40002380
$ws.Range("L99").Value = 6925
$ws.Range("M99").Value = -40000882
$ws.Range("N99").Value = -9921

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 766
$ws.Range("I107").Value = 644.5454999999999
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 644.5454999999999
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 1275.4545
$ws.Range("N107").Value = -4940

$ws.Range("H126").Value = 15389792
$ws.Range("I126").Value = 40002380
$ws.Range("J126").Value = 6925
$ws.Range("K126").Value = 120007140
$ws.Range("L126").Value = 20775
$ws.Range("M126").Value = -120004670
$ws.Range("N126").Value = -25715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 636600.1
$ws.Range("I5").Value = 626.125
$ws.Range("J5").Value = 1027968.75
$ws.Range("K5").Value = 1878.375
$ws.Range("L5").Value = 3083906.25
$ws.Range("M5").Value = -1766.375
$ws.Range("N5").Value = -3084130.25

$ws.Range("H113").Value = 4167266.8
$ws.Range("I113").Value = 626.6667
$ws.Range("J113").Value = 8333907
$ws.Range("K113").Value = 1880.0001
$ws.Range("L113").Value = 25001721
$ws.Range("M113").Value = 289.9999
$ws.Range("N113").Value = -25006061

$ws.Range("H132").Value = 1973.579
$ws.Range("I132").Value = 939.9
$ws.Range("J132").Value = 3122.111
$ws.Range("K132").Value = 8459.1
$ws.Range("L132").Value = 28098.999
$ws.Range("M132").Value = -5929.1
$ws.Range("N132").Value = -33158.999

$ws.Range("H135").Value = 636600.1
$ws.Range("I135").Value = 626.125
$ws.Range("J135").Value = 1027968.75
$ws.Range("K135").Value = 5635.125
$ws.Range("L135").Value = 9251718.75
$ws.Range("M135").Value = -3100.125
$ws.Range("N135").Value = -9256788.75

$ws.Range("H137").Value = 2256
$ws.Range("I137").Value = 1426.6666
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 4279.9998
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = 820.0002000000004
$ws.Range("N137").Value = -20700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 33459.332
$ws.Range("J46").Value = 33459.332
$ws.Range("L46").Value = 33459.332
$ws.Range("N46").Value = -33771.332

$ws.Range("H100").Value = 39199.668
$ws.Range("J100").Value = 39199.668
$ws.Range("L100").Value = 39199.668
$ws.Range("N100").Value = -41363.668

$ws.Range("H106").Value = 34250
$ws.Range("J106").Value = 34250
$ws.Range("L106").Value = 34250
$ws.Range("N106").Value = -36774

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3673.16
$ws.Range("I122").Value = 2068.125
$ws.Range("J122").Value = 6526.5557
$ws.Range("K122").Value = 6204.375
$ws.Range("L122").Value = 19579.6671
$ws.Range("M122").Value = -3754.375
$ws.Range("N122").Value = -24479.6671

$ws.Range("H132").Value = 4788.5713
$ws.Range("I132").Value = 991.8570999999999
$ws.Range("J132").Value = 8585.286
$ws.Range("K132").Value = 2975.5713
$ws.Range("L132").Value = 25755.858
$ws.Range("M132").Value = -445.5712999999996
$ws.Range("N132").Value = -30815.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 44899.668
$ws.Range("J80").Value = 44899.668
$ws.Range("L80").Value = 44899.668
$ws.Range("N80").Value = -46895.668

$ws.Range("H83").Value = 44899.668
$ws.Range("J83").Value = 44899.668
$ws.Range("L83").Value = 134699.004
$ws.Range("N83").Value = -144683.004
